# Lab6_ExpectedResultsTable.xlsx - fix iDecode Er (CBZ) values and add two new rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the two CBZ-related columns (F, G) that had the wrong opcode (10110101 / CBZW)
#     and should use the correct opcode (10110100 / CBZ) ---
$ws.Range("F3").Value = "B4FF FF6B"
$ws.Range("G3").Value = "B400 0109"

$ws.Range("F4").Value = "10110100"
$ws.Range("G4").Value = "10110100"

# read_data1 for the CBZ instructions (F = CBZ X11,-5 reads X11=0 ; G = CBZ X9,8 reads X9=20)
$ws.Range("F15").Value = "0"
$ws.Range("G15").Value = "20"

# --- Append two new explanatory rows under the table ---
$ws.Range("G27").Value = "10110100 1111 1111 1111 1111 0110 1011"
$ws.Range("G28").Value = "10110100 0000 0000 0000 0001 0000 1001"
$ws.Range("G27:G28").Style = $ws.Range("K25").Style

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("H32").Select()
